$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Paragraph about "Suponha que a empresa..." - wording tweaks
#    - "de se evitar desperdicios" -> "de evitar desperdicios"
#    - "com o objetivo de definir" -> "com a finalidade de definir"
# ---------------------------------------------------------------
$r1 = $d.Content
$r1.Find.ClearFormatting()
$r1.Find.Replacement.ClearFormatting()
$f1 = $r1.Find.Execute("de se evitar desperdícios", $true, $false, $false, $false, $false, $true, 0, $false, `
    "de evitar desperdícios", 2)
Write-Host "Step1 (se evitar):" $f1

$r2 = $d.Content
$r2.Find.ClearFormatting()
$r2.Find.Replacement.ClearFormatting()
$f2 = $r2.Find.Execute("investigativa com o objetivo de definir", $true, $false, $false, $false, $false, $true, 0, $false, `
    "investigativa com a finalidade de definir", 2)
Write-Host "Step2 (a finalidade):" $f2

# ---------------------------------------------------------------
# 2) Citation paragraphs (MALWEE / ARMARIO INFANTIL) - strip the
#    literal "<" ">" brackets around the bare URLs and collapse the
#    "Acesso em 16 Jan. 2021." tail back into a single run.
#
#    Order matters:
#      a) merge the ". Acesso em 16 Jan. 2021." tail into one run
#         FIRST, while the "<"/">" runs are still present - this way
#         the merge starts from a plain-formatted run and does not
#         inherit the hyperlink's character style.
#      b) only then remove the leading "<" and the trailing ">"
#         (one citation has a stray space after it: "> ").
#
#    ReplaceAll (last arg = 2) takes care of both citation
#    paragraphs in a single call each.
# ---------------------------------------------------------------

# 2a) Merge ". Acesso em 16 Jan. 2021." into a single run
$ra = $d.Content
$ra.Find.ClearFormatting()
$ra.Find.Replacement.ClearFormatting()
$fa = $ra.Find.Execute(". Acesso em 16 Jan. 2021.", $true, $false, $false, $false, $false, $true, 0, $false, `
    ". Acesso em 16 Jan. 2021.", 2)
Write-Host "Step3a (merge Acesso):" $fa

# 2b) Remove the literal "<" right before each bare URL
$rb = $d.Content
$rb.Find.ClearFormatting()
$rb.Find.Replacement.ClearFormatting()
$fb = $rb.Find.Execute("Disponível em: <", $true, $false, $false, $false, $false, $true, 0, $false, `
    "Disponível em: ", 2)
Write-Host "Step3b (remove <):" $fb

# 2c) Remove the literal ">" right after the hyperlink. One of the two
#     citations has a stray trailing space ("> ") that must go too.
$rc1 = $d.Content
$rc1.Find.ClearFormatting()
$rc1.Find.Replacement.ClearFormatting()
$fc1 = $rc1.Find.Execute("> ", $true, $false, $false, $false, $false, $true, 0, $false, `
    "", 2)
Write-Host "Step3c1 (remove '> '):" $fc1

$rc2 = $d.Content
$rc2.Find.ClearFormatting()
$rc2.Find.Replacement.ClearFormatting()
$fc2 = $rc2.Find.Execute(">", $true, $false, $false, $false, $false, $true, 0, $false, `
    "", 2)
Write-Host "Step3c2 (remove remaining >):" $fc2

$d.Save()
